$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Fix the unicode ellipsis character in cell C9 to plain text
$ws.Range("C9").Value = "Own products, 10Micron, PWI, and others"

# Update the active selection on the sheet
$ws.Range("C10").Select()
